# Update the "grid_cell" (AG) column values on the "solar" sheet.
# The grid-cell association for rows 4-28 was reshuffled to a new
# permutation of CHE_* labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$ws.Range("AG4").Value  = "CHE_24"
$ws.Range("AG5").Value  = "CHE_8"
$ws.Range("AG6").Value  = "CHE_5"
$ws.Range("AG7").Value  = "CHE_7"
$ws.Range("AG8").Value  = "CHE_12"
$ws.Range("AG9").Value  = "CHE_13"
$ws.Range("AG10").Value = "CHE_11"
$ws.Range("AG11").Value = "CHE_15"
$ws.Range("AG12").Value = "CHE_25"
$ws.Range("AG13").Value = "CHE_10"
$ws.Range("AG14").Value = "CHE_22"
$ws.Range("AG15").Value = "CHE_9"
$ws.Range("AG16").Value = "CHE_21"
$ws.Range("AG17").Value = "CHE_4"
$ws.Range("AG18").Value = "CHE_0"
$ws.Range("AG19").Value = "CHE_20"
$ws.Range("AG20").Value = "CHE_1"
$ws.Range("AG21").Value = "CHE_6"
$ws.Range("AG22").Value = "CHE_3"
$ws.Range("AG23").Value = "CHE_17"
$ws.Range("AG24").Value = "CHE_19"
$ws.Range("AG25").Value = "CHE_23"
$ws.Range("AG26").Value = "CHE_14"
$ws.Range("AG27").Value = "CHE_18"
$ws.Range("AG28").Value = "CHE_2"
